$d = $word.ActiveDocument

# Locate the existing "Invincibilità momentanea" bullet in the Power-up
# list (paragraph-exact match, ignoring the trailing paragraph mark).
$para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq "Invincibilità momentanea") {
        $para = $d.Paragraphs($i)
        break
    }
}

# Repurpose that bullet for the new text.
$para.Range.Text = "Bloccare i nemici per n secondi"

# Insert a brand-new list paragraph right after it that recreates the
# original "Invincibilità momentanea" bullet (same list style, inherited
# automatically from the paragraph it follows).
$para.Range.InsertParagraphAfter()
$newPara = $para.Next()
$newPara.Range.Text = "Invincibilità momentanea"

# Word keeps a single "_GoBack" bookmark that tracks the last edit
# location; move it from its old spot (next to "TODO:") into the middle
# of the freshly retyped word, splitting the run into "Invincibilità
# momentan" + bookmark + "ea".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$splitPos = $newPara.Range.Start + "Invincibilità momentan".Length
$bkRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bkRange)
